# Apply the "BurndownChartSampleS3" update:
#  - A handful of daily burndown cells (V7, V8, Y10, X12, Y13, Y14) get
#    filled in with the work completed that day, while J11's previous
#    value is cleared (the task's remaining work was moved to a later day).
#  - The sheet view is scrolled down/right to A4, zoomed to 70%, and the
#    current selection becomes Y15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the daily burndown entries (rows 7-14, "Planilha1") ---
$ws.Range("J11").Value = $null
$ws.Range("V7").Value = 4
$ws.Range("V8").Value = 3
$ws.Range("Y10").Value = 1
$ws.Range("X12").Value = 6
$ws.Range("Y13").Value = 4
$ws.Range("Y14").Value = 2

# Make sure the "Real" burndown row (23) and the chart that reads it
# recompute from the new inputs.
$excel.CalculateFull()

# --- Update the view: scrolled to A4, zoomed to 70%, Y15 selected ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 70
$ws.Range("Y15").Select()
